$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D14").Value = "2016-03-09 04:54:48"
$wsZhCn.Range("G14").Value = "2016-03-09 04:55:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D14").Value = "2016-03-09 04:54:51"
$wsDeDe.Range("G14").Value = "2016-03-09 04:55:47"
